# Add Roll (col B) and Email (new last col) to the absentees sheet, and
# reorder Akshit above Sakshi (per the commit "Added rollnumber and email
# to absentees").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -----------------------------------------------------------
$ws.Cells.Item(1,1).Value = "Name"
$ws.Cells.Item(1,2).Value = "Roll"
$ws.Cells.Item(1,3).Value = "Class"
$ws.Cells.Item(1,4).Value = "Date"
$ws.Cells.Item(1,5).Value = "Email"

# --- Row 2: Nandita --------------------------------------------------------
$ws.Cells.Item(2,1).Value = "Nandita"
$ws.Cells.Item(2,2).NumberFormat = "@"
$ws.Cells.Item(2,2).Value = "1811084"
$ws.Cells.Item(2,3).Value = "MIP"
$ws.Cells.Item(2,4).Value = "Thu Apr 22"
$ws.Cells.Item(2,5).Value = "nandita.kadam@somaiya.edu"

# --- Row 3: Aditya ----------------------------------------------------------
$ws.Cells.Item(3,1).Value = "Aditya"
$ws.Cells.Item(3,2).NumberFormat = "@"
$ws.Cells.Item(3,2).Value = "1811001"
$ws.Cells.Item(3,3).Value = "DSIP"
$ws.Cells.Item(3,4).Value = "Thu Apr 22"
$ws.Cells.Item(3,5).Value = "aditya.pradhan@somaiya.edu"

# --- Row 4: Akshit (moved up from row 5) -----------------------------------
$ws.Cells.Item(4,1).Value = "Akshit"
$ws.Cells.Item(4,2).NumberFormat = "@"
$ws.Cells.Item(4,2).Value = "1811042"
$ws.Cells.Item(4,3).Value = "AI"
$ws.Cells.Item(4,4).Value = "Thu Apr 22"
$ws.Cells.Item(4,5).Value = "akshit.gs@somaiya.edu"

# --- Row 5: Sakshi (moved down from row 4) ---------------------------------
$ws.Cells.Item(5,1).Value = "Sakshi"
$ws.Cells.Item(5,2).NumberFormat = "@"
$ws.Cells.Item(5,2).Value = "1811053"
$ws.Cells.Item(5,3).Value = "AI"
$ws.Cells.Item(5,4).Value = "Thu Apr 22"
$ws.Cells.Item(5,5).Value = "sakshi@somaiya.edu"

# --- Column widths ----------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 9.166666666666666   # -> width 10
$ws.Columns.Item(2).ColumnWidth = 9.166666666666666   # -> width 10
$ws.Columns.Item(3).ColumnWidth = 4.166666666666667   # -> width 5
$ws.Columns.Item(4).ColumnWidth = 14.166666666666666  # -> width 15
$ws.Columns.Item(5).ColumnWidth = 29.166666666666668  # -> width 30
